$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Family")

# Shift the "comments" header from F1 to G1, and add a new "mutations" header in F1.
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = "mutations"

# Widen column F and add the new column G width.
$ws.Columns.Item(6).ColumnWidth = 51.166666666666664
$ws.Columns.Item(7).ColumnWidth = 28.998697916666668

# New diagnosis data for the two other families.
$ws.Range("E3").Value = "Diagnosis 1, Diagnosis 2 test"
$ws.Range("E4").Value = "Diagnosis 2 TEST"
$ws.Range("E4").Font.Color = 0

# New mutations data for the two other families.
$ws.Range("F3").Value = "Gene1 + mutation in Gene1, Gene2 + mutation in Gene2"
$ws.Range("F4").Value = "Gene34 + mutation in Gene34, Gene29 + mutation in Gene29"

# Update the sheet selection/view state.
$ws.Activate()
$ws.Range("G1").Select()

# "clinical values" sheet gains a narrow helper column.
$ws3 = $wb.Worksheets.Item("clinical values")
$ws3.Columns.Item(9).ColumnWidth = 10.166666666666666
